$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - NLH / No-Limit Hold'em variants
$ws.Range("A1").Value = "NLH"
$ws.Range("B1").Value = "NLH"
$ws.Range("C1").Value = "NO-LIMIT HOLD'EM"
$ws.Range("D1").Value = "NL HOLDEM"
$ws.Range("E1").Value = "NL HOLDEM"
$ws.Range("F1").Value = "N/L HOLDEM"
$ws.Range("G1").Value = "NL HOLD'EM"
$ws.Range("H1").Value = "N/L HOLDEM"
$ws.Range("I1").Value = "NO LIMIT HOLDEM"
$ws.Range("J1").Value = "NL Holdem"
$ws.Range("K1").Value = "N/L Holdem"
$ws.Range("L1").Value = "No Limit Holdem"

# Row 2 - LHE / Limit Hold'em variants
$ws.Range("A2").Value = "LHE"
$ws.Range("B2").Value = "LHE"
$ws.Range("C2").Value = "LIMIT HOLD'EM"
$ws.Range("D2").Value = "LIMIT HOLDEM"
$ws.Range("E2").Value = "LIMIT HOLDEM"
$ws.Range("F2").Value = "LIMIT HOLD'EM"
$ws.Range("G2").Value = "Limit Holdem"

# Row 3 - PLO / Pot Limit Omaha
$ws.Range("A3").Value = "PLO"
$ws.Range("B3").Value = "PLO"
$ws.Range("C3").Value = "POT LIMIT OMAHA"

# Row 4 - PLO-H/L / Pot Limit Omaha High/Low
$ws.Range("A4").Value = "PLO-H/L"
$ws.Range("B4").Value = "PLO-H/LOW"
$ws.Range("C4").Value = "POT LIMIT OMAHA HIGH/LOW"

# Row 5 - OMA H/L / Omaha High/Low
$ws.Range("A5").Value = "OMA H/L"
$ws.Range("B5").Value = "OMA/HI-LOW"
$ws.Range("C5").Value = "OMAHA HIGH/LOW"

# Row 6 - LIMIT OMAHA HI / Limit Omaha
$ws.Range("A6").Value = "LIMIT OMAHA HI"
$ws.Range("B6").Value = "LIMIT OMAHA/HI"
$ws.Range("C6").Value = "LIMIT OMAHA"

# Row 7 - PLO/8 / Pot Limit Omaha 8 or Better
$ws.Range("A7").Value = "PLO/8"
$ws.Range("B7").Value = "PLO/8"
$ws.Range("C7").Value = "POT LIMIT OMAHA 8 OR BETTER"

# Row 8 - LO8 / Limit Omaha 8 or Better
$ws.Range("A8").Value = "LO8"
$ws.Range("B8").Value = "LO8"
$ws.Range("C8").Value = "LIMIT OMAHA 8 OR BETTER"

# Resize columns to fit new content (widths mirror the source workbook's best-fit autosize)
$ws.Columns.Item(1).ColumnWidth = 14.0
$ws.Columns.Item(2).ColumnWidth = 14.333333333333334
$ws.Columns.Item(3).ColumnWidth = 27.0
$ws.Columns.Item(4).ColumnWidth = 12.666666666666666
$ws.Columns.Item(5).ColumnWidth = 12.666666666666666
$ws.Columns.Item(6).ColumnWidth = 13.166666666666666
$ws.Columns.Item(7).ColumnWidth = 10.666666666666666
$ws.Columns.Item(8).ColumnWidth = 11.0
$ws.Columns.Item(9).ColumnWidth = 15.833333333333334
$ws.Columns.Item(10).ColumnWidth = 9.166666666666666
$ws.Columns.Item(11).ColumnWidth = 9.166666666666666

# Update selection to match the edited workbook
$ws.Range("L1").Select() | Out-Null
